# Apply edits described by the commit:
# "Started creating the battery module, will work on the payload tomorrow"
#
# Changes:
#  - Inputs!C5: Endurance value 1.2 -> 2
#  - Inputs!B6: Payload/MTOW selector "Payload Weight" -> "MTOW"
#  - Inputs!C6: corresponding value 0.25 -> 20
#  - Inputs sheet selection moves to E15
# Dependent formulas on Inputs and export_ready_inputs sheets recalc automatically.

$wb = $excel.ActiveWorkbook

$wsInputs = $wb.Worksheets.Item("Inputs")
$wsInputs.Activate()

# Update the Endurance input value
$wsInputs.Range("C5").Value = 2

# Switch the Payload/MTOW selector to "MTOW" and update its value
$wsInputs.Range("B6").Value = "MTOW"
$wsInputs.Range("C6").Value = 20

# Move the active selection to E15, matching the saved cursor position
$wsInputs.Range("E15").Select()

$excel.Calculate()

$wb.Save()
